# WorkLog_ns.xlsx update: "add dashboard research and update logs"
#
# Fills in the previously-blank log rows 43-46 (Tuesday block) and 48-53
# (Wednesday block) on the "Week 3" sheet with new time-log entries, and
# leaves the selection parked on C54 (matching where the author's cursor
# ended up after the edit). Dependent SUM() formulas (D46, D58, C94)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 3")

# --- Tuesday block (rows 43-46) --- (column A, then B, then C; top to bottom)
$ws.Range("A43").Value = "Update logs and github"
$ws.Range("B43").Value = "Project Management"
$ws.Range("C43").Value = 0.25

$ws.Range("A44").Value = "Review anomaly detection research PR"
$ws.Range("B44").Value = "Anomaly Detection Research"
$ws.Range("C44").Value = 1

$ws.Range("A45").Value = "Presentation for Tech Safety BC meeting Friday"
$ws.Range("B45").Value = "Project Documents/Client Meetings"
$ws.Range("C45").Value = 1.5

$ws.Range("A46").Value = "Review EDA PR"
$ws.Range("B46").Value = "EDA"
$ws.Range("C46").Value = 0.5

# --- Wednesday block (rows 48-53) ---
# Row 48's Activity Type (B) was filled in before its Description (A).
$ws.Range("B48").Value = "Dashboard Research"
$ws.Range("A48").Value = "Initial research on Grafana to understand dashboard level-of-effort"
$ws.Range("C48").Value = 4

$ws.Range("A49").Value = "Stand-up"
$ws.Range("B49").Value = "Internal Meetings"
$ws.Range("C49").Value = 0.5

$ws.Range("A50").Value = "Sprint planning meeting"
$ws.Range("B50").Value = "Client Meetings"
$ws.Range("C50").Value = 1

$ws.Range("A51").Value = "Helping UDL with data streaming parsing"
$ws.Range("B51").Value = "Streaming Parsing support for UDL"
$ws.Range("C51").Value = 1.25

# Row 53 was filled in before row 52.
$ws.Range("A53").Value = "Document Grafana research (sufficient understanding for now)"
$ws.Range("B53").Value = "Dashboard Research"
$ws.Range("C53").Value = 1.25

$ws.Range("A52").Value = "Discuss anomaly detection approach w/ Ryan and generate schematic"
$ws.Range("B52").Value = "Anomaly Detection Research"
$ws.Range("C52").Value = 1.25

# Leave the cursor where the author's saved selection ended up.
$ws.Activate()
$ws.Range("C54").Select()
